$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 285 (「惑星地球はあなた方を歓迎します」...),
# shifting all subsequent rows up by one.
$ws.Rows.Item(285).Delete()
